$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.674.32'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.845.23'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.10%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.40'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('E6').Value = '  +0.09%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4318'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +0.26%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3710'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +1.77%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07344'
$ws.Range('D9').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8804'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +0.31%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.00'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '1.846.03'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('E13').Value = '  +2.85%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.613'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +1.28%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06955'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +0.66%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +0.09%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.21'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +1.49%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009063'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +0.53%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.10%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.60'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('D21').Value = '27.701.84'
$ws.Range('E21').Value = '  +0.09%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.138'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +3.56%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +5.98%  '
$ws.Range('D24').Value = '2.078.99'
$ws.Range('E24').Value = '  -0.43%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.988'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.20%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.92'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('E27').Value = '  +0.80%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.316'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +0.45%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.37'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -4.65%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.877'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +1.31%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08936'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7881'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +2.91%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.617'
$ws.Range('D33').Style = $origStyle
$ws.Range('E34').Value = '  +6.51%  '
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  +0.18%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05444'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('E38').Value = '  +1.19%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01968'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +1.39%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.843'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +0.55%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5182'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +1.79%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1692'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +2.34%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.788'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +0.38%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.640'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('E45').Value = '  +2.87%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4795'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +2.42%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '106.86'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +1.78%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06559'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  +0.14%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.667'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +2.53%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.839'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +5.11%  '
